# Reorganize the Rieslings list: replace the wine names / sugar / volume /
# price data in rows 2-20 with the new data set, and remove the now-unused
# rows 21-25 (the list shrank from 24 wines to 19 wines).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  Name='Joseph Cattin Riesling 2014';                                           B=5;   C=750; D=14.95}
    @{Row=3;  Name='Schieferkopf Lieu-Dit Fels Riesling 2011';                               B=5;   C=750; D=49}
    @{Row=4;  Name='Gustave Lorentz Réserve Riesling 2014';                                  B=6;   C=750; D=18.95}
    @{Row=5;  Name='Boeckel Wiebelsberg Riesling 2012';                                      B=6;   C=750; D=35}
    @{Row=6;  Name='Baron de Hoen Réserve Riesling 2014';                                    B=7;   C=750; D=14.95}
    @{Row=7;  Name='Cave de Beblenheim Heimberger Réserve Riesling';                         B=7;   C=750; D=15.95}
    @{Row=8;  Name='Gustave Lorentz Riesling Cuvee Amethyste';                               B=7;   C=750; D=16.3}
    @{Row=9;  Name='Henri Ehrhart Réserve Particulière Riesling 2013';                       B=7;   C=750; D=16.95}
    @{Row=10; Name='Trimbach Riesling 2012';                                                 B=7;   C=750; D=21.95}
    @{Row=11; Name='Boeckel Brandluft Riesling 2012';                                        B=8;   C=750; D=14.75}
    @{Row=12; Name='Jean Geiler Réserve Particulière Riesling 2013';                         B=8;   C=750; D=14.95}
    @{Row=13; Name='Willm Réserve Riesling';                                                 B=8;   C=750; D=15.95}
    @{Row=14; Name='J. Fritsch Riesling 2014';                                               B=8;   C=750; D=17.25}
    @{Row=15; Name='Hugel & Fils Jubilee Riesling 2009';                                     B=8;   C=750; D=55}
    @{Row=16; Name='Pierre Sparr Lieu Dit Altenbourg Riesling 2013';                         B=9;   C=750; D=16.95}
    @{Row=17; Name='Koenig Riesling Kp M 2014';                                              B=10;  C=750; D=18.95}
    @{Row=18; Name='Pierre Sparr Schoenenbourg Riesling 2011';                               B=11;  C=750; D=22.95}
    @{Row=19; Name='Domaine Pfister Silberberg Sélection de Grains Nobles Riesling 2007';    B=100; C=500; D=56.25}
    @{Row=20; Name='Clos St. Landelin Vorbourg Sélection de Grains Nobles Riesling 2007';    B=180; C=500; D=51.75}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value2 = $item.Name
    $ws.Cells.Item($r, 2).Value2 = $item.B
    $ws.Cells.Item($r, 3).Value2 = $item.C
    $ws.Cells.Item($r, 4).Value2 = $item.D
}

# Rows 21-25 no longer correspond to any wine in the reorganized list.
$ws.Rows("21:25").Delete()
